$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BF holds a game-date label (row 1 header "Date") followed by date
# strings for each team row (rows 2-31). Those values were originally
# mis-derived as "5-15-2011-12" (file name mashed together) instead of the
# actual game date. Correct them to the ISO-style date string "2012-05-15".
#
# NOTE: the source values are plain text, not real Excel dates, so we force
# the cell to Text format before assigning the new value (otherwise Excel's
# automatic date recognition would turn the literal string into a date
# serial number) and then restore the cell's style afterwards so no visible
# formatting change is introduced.
$oldValue = "5-15-2011-12"
$newValue = "2012-05-15"

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 58)
    if ($cell.Value() -eq $oldValue) {
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
        $cell.Style = "Normal"
    }
}
